$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("B2").Value = "prodOfdj"
$ws.Range("B3").Value = "prodhOje"
$ws.Range("B5").Value = "prodGldp"
